# Fruta / hortaliza, semanal
# Rows 4-9 hold one weekly price record each (Vega Monumental Concepción - Tuna).
# The commit reorders which week's figures sit in which row: row 4 swaps
# with row 9, row 5 swaps with row 7, and row 6 swaps with row 8.
# Only the columns that vary per-record change: D (Fecha), L (Calidad),
# M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), Q (Unidad de comercialización),
# S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "Q", "S")
$swaps = @(@(4, 9), @(5, 7), @(6, 8))

foreach ($pair in $swaps) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        # Value2 is used for the read because it reliably returns the raw
        # cell content (number or string); Value is kept for the write
        # since both forms update the cell correctly.
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
